$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 83 ("Group Discussion" filler row) - everything below shifts up by one.
# Old row 84 ("Final Assignment") becomes the new row 83.
$ws.Rows(83).Delete()

# Fill in the previously-blank Actual Start / Actual Finish / Actual Sign-off dates
# for the Week 8 tasks (rows 72-83 after the shift).
$ws.Range("I72").Value = 44685
$ws.Range("K72").Value = 44684

$ws.Range("I73").Value = 44685
$ws.Range("K73").Value = 44684

$ws.Range("K74").Value = 44684

$ws.Range("K75").Value = 44685

$ws.Range("H76").Value = 44684
$ws.Range("I76").Value = 44684
$ws.Range("K76").Value = 44684

$ws.Range("I77").Value = 44684
$ws.Range("K77").Value = 44684

$ws.Range("H78").Value = 44684
$ws.Range("I78").Value = 44684
$ws.Range("K78").Value = 44684

$ws.Range("H79").Value = 44685
$ws.Range("I79").Value = 44685
$ws.Range("K79").Value = 44685

$ws.Range("H80").Value = 44684
$ws.Range("I80").Value = 44684
$ws.Range("K80").Value = 44684

$ws.Range("H81").Value = 44684
$ws.Range("I81").Value = 44684
$ws.Range("K81").Value = 44684

$ws.Range("H82").Value = 44685
$ws.Range("I82").Value = 44685
$ws.Range("K82").Value = 44685

# Row 83 is the old row 84 ("Final Assignment") shifted up.
$ws.Range("H83").Value = 44685
$ws.Range("I83").Value = 44685
$ws.Range("K83").Value = 44685

# Re-point the AutoFilter / FilterDatabase range now that the sheet is one row shorter.
$ws.AutoFilterMode = $false
$ws.Range("D5:E83").AutoFilter()

foreach ($n in $wb.Names) {
    if ($n.Name -eq "Project_Plan!_FilterDatabase") {
        $n.RefersTo = "=Project_Plan!`$D`$5:`$E`$83"
    }
}

# Re-point the conditional formatting range (it still targets the old M6:AQ84).
$oldCf = $ws.Range("M6:AQ84").FormatConditions
for ($i = 1; $i -le $oldCf.Count; $i++) {
    $oldCf.Item($i).ModifyAppliesToRange($ws.Range("M6:AQ83"))
}
